$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.339747122270978
$ws.Range("C2").Value = 0.3051189489374053
$ws.Range("D2").Value = 0.01080933977811327
$ws.Range("F2").Value = 4.512974865459626
$ws.Range("G2").Value = 0.002633206879539633
$ws.Range("J2").Value = 0.155162212574913
$ws.Range("L2").Value = 0.4533808575927338
$ws.Range("N2").Value = 2.657753170492086

$ws.Range("B3").Value = 4.215118333504734
$ws.Range("C3").Value = 0.2778826528646903
$ws.Range("D3").Value = 0.009961498735012952
$ws.Range("F3").Value = 4.482044098002049
$ws.Range("G3").Value = 0.002639782201115738
$ws.Range("J3").Value = 0.1559184989454039
$ws.Range("L3").Value = 0.4484557097681545
$ws.Range("N3").Value = 2.674764486199436

$ws.Range("B4").Value = 4.14125960592969
$ws.Range("C4").Value = 0.2613369118196829
$ws.Range("D4").Value = 0.009437064542769491
$ws.Range("F4").Value = 4.465391231654564
$ws.Range("G4").Value = 0.002644030233316325
$ws.Range("J4").Value = 0.156418857809923
$ws.Range("L4").Value = 0.4456638698741386
$ws.Range("N4").Value = 2.685984874762099

$ws.Range("B5").Value = 4.111828701034426
$ws.Range("C5").Value = 0.254638271367611
$ws.Range("D5").Value = 0.009222315762542621
$ws.Range("F5").Value = 4.459190876041291
$ws.Range("G5").Value = 0.002645814530827655
$ws.Range("J5").Value = 0.1566318177609443
$ws.Range("L5").Value = 0.4445845032795859
$ws.Range("N5").Value = 2.690751636343457

$ws.Range("B6").Value = 4.106981936425143
$ws.Range("C6").Value = 0.2535285962989917
$ws.Range("D6").Value = 0.009186592248781267
$ws.Range("F6").Value = 4.458196628190152
$ws.Range("G6").Value = 0.002646114030199851
$ws.Range("J6").Value = 0.1566677270229704
$ws.Range("L6").Value = 0.4444087965739385
$ws.Range("N6").Value = 2.691554874700195

$ws.Range("B7").Value = 4.140859992810647
$ws.Range("C7").Value = 0.2612463948972277
$ws.Range("D7").Value = 0.009434172657545759
$ws.Range("F7").Value = 4.465305242477328
$ws.Range("G7").Value = 0.002644054081312151
$ws.Range("J7").Value = 0.1564216931692801
$ws.Range("L7").Value = 0.4456490770573396
$ws.Range("N7").Value = 2.686048374806418

$ws.Range("B8").Value = 4.296220525678393
$ws.Range("C8").Value = 0.2956905868352067
$ws.Range("D8").Value = 0.01051776339519961
$ws.Range("F8").Value = 4.501823129273703
$ws.Range("G8").Value = 0.00263543042400642
$ws.Range("J8").Value = 0.155415514026835
$ws.Range("L8").Value = 0.4516344273277042
$ws.Range("N8").Value = 2.663457374899551

$ws.Range("B9").Value = 4.622165468348669
$ws.Range("C9").Value = 0.3646808529389887
$ws.Range("D9").Value = 0.01261522641012292
$ws.Range("F9").Value = 4.592104397121801
$ws.Range("G9").Value = 0.002620182912156981
$ws.Range("J9").Value = 0.1537276559673213
$ws.Range("L9").Value = 0.4652184771554744
$ws.Range("N9").Value = 2.625336144648699

$ws.Range("B10").Value = 4.874850676987307
$ws.Range("C10").Value = 0.4163064428419148
$ws.Range("D10").Value = 0.01414396595948375
$ws.Range("F10").Value = 4.669986841244821
$ws.Range("G10").Value = 0.002609982285707182
$ws.Range("J10").Value = 0.1526610138409552
$ws.Range("L10").Value = 0.4763324052972422
$ws.Range("N10").Value = 2.601133527402794

$ws.Range("B11").Value = 4.992727792651976
$ws.Range("C11").Value = 0.4400089436870758
$ws.Range("D11").Value = 0.01483774719773834
$ws.Range("F11").Value = 4.707964771814773
$ws.Range("G11").Value = 0.002605556623828309
$ws.Range("J11").Value = 0.152213334038251
$ws.Range("L11").Value = 0.4816365027267295
$ws.Range("N11").Value = 2.590957686691596

$ws.Range("B12").Value = 5.037789968536686
$ws.Range("C12").Value = 0.4490167354611572
$ws.Range("D12").Value = 0.01510030567352771
$ws.Range("F12").Value = 4.722715513351773
$ws.Range("G12").Value = 0.002603911404922317
$ws.Range("J12").Value = 0.1520492010771441
$ws.Range("L12").Value = 0.4836808618936175
$ws.Range("N12").Value = 2.587225020842112

$ws.Range("B13").Value = 5.028066079615996
$ws.Range("C13").Value = 0.447075301805512
$ws.Range("D13").Value = 0.01504376496867366
$ws.Range("F13").Value = 4.719522206139317
$ws.Range("G13").Value = 0.002604264370342166
$ws.Range("J13").Value = 0.1520843102339988
$ws.Range("L13").Value = 0.48323897806074
$ws.Range("N13").Value = 2.588023536477095

$ws.Range("B14").Value = 4.996426553934384
$ws.Range("C14").Value = 0.4407493704611056
$ws.Range("D14").Value = 0.01485935083895384
$ws.Range("F14").Value = 4.709170906123063
$ws.Range("G14").Value = 0.002605420657026482
$ws.Range("J14").Value = 0.1521997226491791
$ws.Range("L14").Value = 0.481803975025997
$ws.Range("N14").Value = 2.590648174487242

$ws.Range("B15").Value = 4.977101839108627
$ws.Range("C15").Value = 0.436878769898442
$ws.Range("D15").Value = 0.0147463731055808
$ws.Range("F15").Value = 4.702878617088743
$ws.Range("G15").Value = 0.002606132906275785
$ws.Range("J15").Value = 0.1522711183652419
$ws.Range("L15").Value = 0.4809296613042306
$ws.Range("N15").Value = 2.59227158333016

$ws.Range("B16").Value = 4.867206367503115
$ws.Range("C16").Value = 0.4147618833133606
$ws.Range("D16").Value = 0.01409859898255661
$ws.Range("F16").Value = 4.667556421461057
$ws.Range("G16").Value = 0.002610275816568339
$ws.Range("J16").Value = 0.1526910257832945
$ws.Range("L16").Value = 0.4759907756574222
$ws.Range("N16").Value = 2.601815391172607

$ws.Range("B17").Value = 4.800541743449287
$ws.Range("C17").Value = 0.401250266902764
$ws.Range("D17").Value = 0.01370084158028817
$ws.Range("F17").Value = 4.646542190097819
$ws.Range("G17").Value = 0.002612872205010075
$ws.Range("J17").Value = 0.1529582371233928
$ws.Range("L17").Value = 0.4730246018813773
$ws.Range("N17").Value = 2.60788437143286

$ws.Range("B18").Value = 4.762473527613167
$ws.Range("C18").Value = 0.3934991870504518
$ws.Range("D18").Value = 0.01347190213691363
$ws.Range("F18").Value = 4.634695095344199
$ws.Range("G18").Value = 0.002614385794217561
$ws.Range("J18").Value = 0.1531154638254009
$ws.Range("L18").Value = 0.471341903205186
$ws.Range("N18").Value = 2.611453546490409

$ws.Range("B19").Value = 4.749631500347277
$ws.Range("C19").Value = 0.390878290178307
$ws.Range("D19").Value = 0.01339435741539319
$ws.Range("F19").Value = 4.630724962783006
$ws.Range("G19").Value = 0.002614901747092363
$ws.Range("J19").Value = 0.1531693051667293
$ws.Range("L19").Value = 0.4707761803741306
$ws.Range("N19").Value = 2.612675459866836

$ws.Range("B20").Value = 4.807609767411066
$ws.Range("C20").Value = 0.4026864794616927
$ws.Range("D20").Value = 0.01374319953868053
$ws.Range("F20").Value = 4.648754359067027
$ws.Range("G20").Value = 0.002612593724118724
$ws.Range("J20").Value = 0.152929426319929
$ws.Range("L20").Value = 0.4733379371054127
$ws.Range("N20").Value = 2.607230192482305

$ws.Range("B21").Value = 5.005708296285775
$ws.Range("C21").Value = 0.4426065698791604
$ws.Range("D21").Value = 0.01491352150971892
$ws.Range("F21").Value = 4.712201287860182
$ws.Range("G21").Value = 0.002605080196458224
$ws.Range("J21").Value = 0.1521656768789938
$ws.Range("L21").Value = 0.4822244973931475
$ws.Range("N21").Value = 2.589873972611599

$ws.Range("B22").Value = 5.137654105490697
$ws.Range("C22").Value = 0.4688845719797428
$ws.Range("D22").Value = 0.01567748960534487
$ws.Range("F22").Value = 4.755821568043757
$ws.Range("G22").Value = 0.002600348429041555
$ws.Range("J22").Value = 0.1516979588872793
$ws.Range("L22").Value = 0.4882411742654824
$ws.Range("N22").Value = 2.579234602006068

$ws.Range("B23").Value = 5.067004340597236
$ws.Range("C23").Value = 0.4548420264180777
$ws.Range("D23").Value = 0.01526980293537861
$ws.Range("F23").Value = 4.73234256555034
$ws.Range("G23").Value = 0.00260285756663757
$ws.Range("J23").Value = 0.1519447138275503
$ws.Range("L23").Value = 0.4850108195236515
$ws.Range("N23").Value = 2.584848367970238

$ws.Range("B24").Value = 4.804413509911228
$ws.Range("C24").Value = 0.4020371151000859
$ws.Range("D24").Value = 0.01372405033163204
$ws.Range("F24").Value = 4.647753507998772
$ws.Range("G24").Value = 0.002612719560271
$ws.Range("J24").Value = 0.1529424404599169
$ws.Range("L24").Value = 0.4731962078825944
$ws.Range("N24").Value = 2.607525697544986

$ws.Range("B25").Value = 4.531685867951182
$ws.Range("C25").Value = 0.3458565563292382
$ws.Range("D25").Value = 0.01205028749491177
$ws.Range("F25").Value = 4.565664692492135
$ws.Range("G25").Value = 0.002624130964374844
$ws.Range("J25").Value = 0.1541537774416106
$ws.Range("L25").Value = 0.4613451254826231
$ws.Range("N25").Value = 2.634983513626025
